# Update "想去人数" (interest count) figures for a handful of events.
# These values are duplicated between the "展览" sheet and the "全部类型"
# sheet (which aggregates rows from all the other sheets), so both need
# to be updated to keep the workbook consistent.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13231
$ws1.Range("F7").Value = 111
$ws1.Range("F11").Value = 13170
$ws1.Range("F27").Value = 84

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13231
$ws4.Range("F8").Value = 111
$ws4.Range("F12").Value = 13170
$ws4.Range("F30").Value = 84
